# Tendencias Climáticas y Predicción Global de Anomalías de Temperatura
# "Se realizo la EDA (Exploratory Data Analysis)"
#
# The FASE 1 heading is updated from a "under construction" marker (🛠️)
# to a "done" checkmark (✅), which reflows the run/proofErr structure
# the way Word's own editor does when the paragraph is retyped:
#   🛠️ FASE 1: Preparación del entorno profesional
#     ->
#   ✅  FASE 1: Preparación del entorno profesional
# (note the doubled space before "FASE" and the wrapping gramStart/gramEnd
#  proofing marks Word leaves behind around "FASE").

$d = $word.ActiveDocument

# Locate the "FASE 1" heading paragraph by its (still-unique) text so the
# edit does not depend on a brittle, hard-coded paragraph index.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*FASE 1:*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'FASE 1' heading paragraph"
}

# Rebuild the paragraph's contents exactly as Word leaves them after the
# emoji swap: a gramStart mark, the new checkmark run (no longer bold),
# a lone-space run, the (still bold) " FASE" run, a gramEnd mark, and the
# remaining (still bold) " 1: Preparación del entorno profesional" run.
$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body><w:p w14:paraId="154035E6" w14:textId="77777777" w:rsidR="00B45267" w:rsidRPr="00B45267" w:rsidRDefault="00B45267" w:rsidP="00B45267"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00B45267"><w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></w:rPr><w:t>&#x2705;</w:t></w:r><w:r w:rsidRPr="00B45267"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00B45267"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> FASE</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00B45267"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> 1: Preparaci&#xF3;n del entorno profesional</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($newParaXml)
